$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for Tyler Anderson right after the header ---
$ws.Rows("2:2").Insert()

# --- Row 2: Tyler Anderson (new), written left to right ---
$ws.Range("A2").Value = "Tyler Anderson"
$ws.Range("B2").Value = "anderty01"
$ws.Range("C2").Value = "June 03 2017"
$ws.Range("D2").Value = "Knee"
$ws.Range("E2").Value = "Anderson has been placed on the 10-day disabled list with knee inflammation it is unknown when he will rejoin the team."

# --- Chad Bettis row: refresh the Last.Updated date ---
$ws.Range("C3").Value = "June 03 2017"

# --- David Dahl row: date, then detail text ---
$ws.Range("C4").Value = "May 29 2017"
$ws.Range("E4").Value = "Dahl is on the 10-day disabled list with a stress reaction of his sixth rib and is likely to remain sidelined until sometime in June."

# --- Jon Gray row: date, then detail text ---
$ws.Range("C5").Value = "June 01 2017"
$ws.Range("E5").Value = "Gray has been placed on the 10-day disabled list with a stress fracture in his left foot and is expected to be sidelined until the end of June."

# --- Tom Murphy row: date only for now ---
$ws.Range("C6").Value = "May 30 2017"

# --- Row 7: Adam Ottavino (new) ---
$ws.Range("A7").Value = "Adam Ottavino"
$ws.Range("B7").Value = "ottavad01"
$ws.Range("C7").Value = "May 30 2017"
$ws.Range("D7").Value = "Shoulder"

# --- Detail text updates (Bettis, Murphy, Ottavino) ---
$ws.Range("E3").Value = "Bettis is on the 60-day disabled list while recovering from testicular cancer and it is unknown as to when `nhe will be ready to rejoin the team."
$ws.Range("E6").Value = "Murphy is on the 10-day disabled list while he recovers from a hairline fracture in his wrist but is expected to `nreturn by the end of June."
$ws.Range("E7").Value = "Ottavino has been placed on the 10-day disabled list with a inflammation in his right shoulder and is without `na timetable for return."

# --- Wrap text + row heights for the long description cells ---
$ws.Range("E2").WrapText = $true
$ws.Rows("2:2").RowHeight = 30

$ws.Range("E4").WrapText = $true
$ws.Rows("4:4").RowHeight = 30

$ws.Range("E6").WrapText = $true
$ws.Rows("6:6").RowHeight = 45

$ws.Range("E7").WrapText = $true
$ws.Rows("7:7").RowHeight = 45

# --- Header row: turn off wrap on the Injury.Details header cell ---
$ws.Range("E1").WrapText = $false

# --- Column B width (PlayerID column) ---
$ws.Columns("B:B").ColumnWidth = 9.75

# --- Page orientation ---
$ws.PageSetup.Orientation = 1

# --- Reset window scroll position / selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("E5").Select()
